$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the price/volume columns so numeric-looking strings
# (e.g. "1.003", "0.3955") are preserved as text instead of being parsed as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '24.556.06'
$ws.Range('E2').Value = '  -1.21%  '
$ws.Range('D3').Value = '1.673.04'
$ws.Range('E3').Value = '  -1.99%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '314.08'
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('D7').Value = '0.3955'
$ws.Range('E7').Value = '  -1.58%  '
$ws.Range('D8').Value = '0.3936'
$ws.Range('D9').Value = '1.003'
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('D10').Value = '1.394'
$ws.Range('E10').Value = '  -5.16%  '
$ws.Range('D11').Value = '50.30'
$ws.Range('E11').Value = '  -6.41%  '
$ws.Range('D12').Value = '0.08637'
$ws.Range('E12').Value = '  -1.80%  '
$ws.Range('D13').Value = '25.32'
$ws.Range('E13').Value = '  -3.85%  '
$ws.Range('D14').Value = '7.296'
$ws.Range('E14').Value = '  -2.77%  '
$ws.Range('D15').Value = '0.00001314'
$ws.Range('E15').Value = '  -2.10%  '
$ws.Range('D16').Value = '7.658'
$ws.Range('E16').Value = '  -4.32%  '
$ws.Range('D17').Value = '1.676.13'
$ws.Range('E17').Value = '  +3.43%  '
$ws.Range('D18').Value = '93.90'
$ws.Range('E18').Value = '  -1.65%  '
$ws.Range('D19').Value = '0.07010'
$ws.Range('E19').Value = '  -2.38%  '
$ws.Range('D20').Value = '21.19'
$ws.Range('E20').Value = '  +1.12%  '
$ws.Range('D21').Value = '7.065'
$ws.Range('E21').Value = '  -3.06%  '
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('E23').Value = '  -3.98%  '
$ws.Range('D24').Value = '24.561.48'
$ws.Range('E24').Value = '  -1.17%  '
$ws.Range('D25').Value = '2.346'
$ws.Range('E25').Value = '  +0.33%  '
$ws.Range('D26').Value = '2.759'
$ws.Range('E26').Value = '  -4.51%  '
$ws.Range('D27').Value = '23.00'
$ws.Range('E27').Value = '  -0.33%  '
$ws.Range('D28').Value = '5.834'
$ws.Range('E28').Value = '  -9.01%  '
$ws.Range('D29').Value = '158.92'
$ws.Range('E29').Value = '  -1.68%  '
$ws.Range('D30').Value = '145.42'
$ws.Range('D31').Value = '8.313'
$ws.Range('E31').Value = '  -0.29%  '
$ws.Range('D32').Value = '2.533'
$ws.Range('E32').Value = '  +10.65%  '
$ws.Range('D33').Value = '1.850.82'
$ws.Range('E33').Value = '  +2.15%  '
$ws.Range('D34').Value = '0.03075'
$ws.Range('E34').Value = '  -3.45%  '
$ws.Range('D35').Value = '0.08251'
$ws.Range('E35').Value = '  -5.16%  '
$ws.Range('D36').Value = '6.902'
$ws.Range('E36').Value = '  -4.28%  '
$ws.Range('D37').Value = '0.2799'
$ws.Range('E37').Value = '  -2.30%  '
$ws.Range('D38').Value = '0.9906'
$ws.Range('E38').Value = '  -3.59%  '
$ws.Range('D39').Value = '0.09634'
$ws.Range('E39').Value = '  +2.11%  '
$ws.Range('D40').Value = '1.511'
$ws.Range('E40').Value = '  +2.01%  '
$ws.Range('D41').Value = '10.28'
$ws.Range('E41').Value = '  -5.04%  '
$ws.Range('D42').Value = '0.7859'
$ws.Range('E42').Value = '  -6.60%  '
$ws.Range('D43').Value = '13.49'
$ws.Range('E43').Value = '  -5.12%  '
$ws.Range('D44').Value = '16.53'
$ws.Range('E44').Value = '  -5.21%  '
$ws.Range('D45').Value = '2.557'
$ws.Range('E45').Value = '  -6.07%  '
$ws.Range('D46').Value = '0.7076'
$ws.Range('E46').Value = '  -4.69%  '
$ws.Range('D47').Value = '4.165'
$ws.Range('E47').Value = '  -1.42%  '
$ws.Range('D48').Value = '0.08632'
$ws.Range('E48').Value = '  +2.90%  '
$ws.Range('E49').Value = '  +0.20%  '
$ws.Range('E50').Value = '  -3.29%  '
$ws.Range('D51').Value = '137.61'
$ws.Range('E51').Value = '  -2.04%  '
